$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D
$ws.Range("D1").Value = "Resistance-"

# Add formulas for D2:D8 referencing column B
$ws.Range("D2").Formula = "=B2-0.000428"
$ws.Range("D3:D8").Formula = "=B3-0.000428"

# Update selection to match target state
$ws.Range("D2:D8").Select()
